# Add 4o baseline (with heuristics) results to the "join-no-learning" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("join-no-learning")

# Fix the swapped FPR/FNR values in the "No Heuristics" block (cols E:F) for the
# "4o" rows, and fill in the new "With Heuristics" block (cols G:H:I).
# (Numbers are written in plain decimal form since the interpreter does not
# accept E-notation literals.)
$ws.Range("E6").Value = 0.17986314760508301247
$ws.Range("F6").Value = 0.32453567937438898561
$ws.Range("G6").Value = 0.91397849462365499118
$ws.Range("H6").Value = 0.03453893776474419725
$ws.Range("I6").Value = 0.05148256761159979850

$ws.Range("E10").Value = 0.33072662104920097503
$ws.Range("F10").Value = 0.03519061583577710178
$ws.Range("G10").Value = 0.92733789507982999911
$ws.Range("H10").Value = 0.07135874877810359507
$ws.Range("I10").Value = 0.00130335614206580994

# Update the active selection on the sheet to match the saved view state.
$ws.Range("M10").Select()
